# Updated symbol list on Sun Jan  1 10:36:03 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (D) and "Volume(1h)" (E) columns on the active sheet
# with the latest scraped values. These columns are stored as plain text
# (not numbers/percentages) in the source data, so each value is written
# with a leading apostrophe to force Excel to keep it as text instead of
# auto-converting it to a Number/Percentage cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'244.13"
$ws.Range("E2").Value  = "'-0.72%"

$ws.Range("D3").Value  = "'27.07"
$ws.Range("E3").Value  = "'3.46%"

$ws.Range("D4").Value  = "'5.156"
$ws.Range("E4").Value  = "'1.16%"

$ws.Range("D5").Value  = "'0.05622"
$ws.Range("E5").Value  = "'0.52%"

$ws.Range("D6").Value  = "'6.470"
$ws.Range("E6").Value  = "'-0.27%"

$ws.Range("D7").Value  = "'0.8175"
$ws.Range("E7").Value  = "'0.31%"

$ws.Range("D8").Value  = "'0.8318"
$ws.Range("E8").Value  = "'-1.89%"

$ws.Range("D9").Value  = "'0.1329"
$ws.Range("E9").Value  = "'-1.00%"

$ws.Range("D10").Value = "'0.06914"
$ws.Range("E10").Value = "'-0.55%"

$ws.Range("D11").Value = "'0.02895"
$ws.Range("E11").Value = "'1.76%"

$ws.Range("D12").Value = "'0.09385"
$ws.Range("E12").Value = "'-0.07%"

$ws.Range("D13").Value = "'0.001523"
$ws.Range("E13").Value = "'0.50%"

$ws.Range("D14").Value = "'0.04288"
$ws.Range("E14").Value = "'-8.56%"

$ws.Range("D15").Value = "'0.0005947"
$ws.Range("E15").Value = "'-93.92%"

$ws.Range("D16").Value = "'0.006110"
$ws.Range("E16").Value = "'-1.30%"

$ws.Range("E17").Value = "'1.58%"

$ws.Range("D18").Value = "'3.021"
$ws.Range("E18").Value = "'0.03%"

$ws.Range("D19").Value = "'2.308"
$ws.Range("E19").Value = "'8.95%"

$ws.Range("D21").Value = "'0.03119"
$ws.Range("E21").Value = "'-2.82%"

$ws.Range("E22").Value = "'-2.13%"

$ws.Range("D23").Value = "'3.736"
$ws.Range("E23").Value = "'-0.23%"

$ws.Range("D25").Value = "'0.001224"
$ws.Range("E25").Value = "'-1.80%"

$ws.Range("D26").Value = "'0.004482"
$ws.Range("E26").Value = "'-2.86%"

$ws.Range("D27").Value = "'0.00009794"

$ws.Range("E28").Value = "'-0.44%"

$ws.Range("D40").Value = "'0.03651"
$ws.Range("E40").Value = "'-0.08%"

$ws.Range("D41").Value = "'0.006057"
$ws.Range("E41").Value = "'77.66%"

$ws.Range("D42").Value = "'0.1052"
$ws.Range("E42").Value = "'-22.95%"

$ws.Range("D43").Value = "'0.002613"
$ws.Range("E43").Value = "'0.75%"

$ws.Range("D44").Value = "'0.008156"
$ws.Range("E44").Value = "'5.06%"

$ws.Range("D45").Value = "'0.00005322"
$ws.Range("E45").Value = "'0.28%"

$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'0.00%"

$ws.Range("E47").Value = "'-35.88%"

$ws.Range("D48").Value = "'0.002635"
$ws.Range("E48").Value = "'28.69%"

$ws.Range("E49").Value = "'0.00%"

$ws.Range("E50").Value = "'0.00%"
